$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3954.7334
$ws.Range("I62").Value = 4126.8335
$ws.Range("J62").Value = 3266.3333
$ws.Range("K62").Value = 4126.8335
$ws.Range("L62").Value = 3266.3333
$ws.Range("M62").Value = -3502.8335
$ws.Range("N62").Value = -4514.3333

$ws.Range("H65").Value = 3954.7334
$ws.Range("I65").Value = 4126.8335
$ws.Range("J65").Value = 3266.3333
$ws.Range("K65").Value = 20634.1675
$ws.Range("L65").Value = 16331.6665
$ws.Range("M65").Value = -17514.1675
$ws.Range("N65").Value = -22571.6665

$ws.Range("H69").Value = 28824.408
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 28824.408
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 86473.224
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -88221.224

$ws.Range("H72").Value = 28824.408
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 28824.408
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 259419.672
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -268155.672

$ws.Range("H87").Value = 19999.953
$ws.Range("J87").Value = 19999.953
$ws.Range("L87").Value = 19999.953
$ws.Range("N87").Value = -22495.953

$ws.Range("H90").Value = 19999.953
$ws.Range("J90").Value = 19999.953
$ws.Range("L90").Value = 59999.859
$ws.Range("N90").Value = -72479.859

$ws.Range("H92").Value = 1333.5652
$ws.Range("I92").Value = 1108.2858
$ws.Range("K92").Value = 1108.2858
$ws.Range("M92").Value = 139.7141999999999

$ws.Range("H100").Value = 5162.9
$ws.Range("I100").Value = 2511.111
$ws.Range("K100").Value = 2511.111
$ws.Range("M100").Value = -1970.111

$ws.Range("H116").Value = 3489.1
$ws.Range("J116").Value = 4347.5
$ws.Range("L116").Value = 4347.5
$ws.Range("N116").Value = -11231.5

$ws.Range("H132").Value = 12932.619
$ws.Range("I132").Value = 1018.4474
$ws.Range("K132").Value = 3055.3422
$ws.Range("M132").Value = -525.3422

$ws.Range("H138").Value = 2403.65
$ws.Range("J138").Value = 3503.1365
$ws.Range("L138").Value = 10509.4095
$ws.Range("N138").Value = -20789.4095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 15734.046
$ws.Range("I2").Value = 25421.154
$ws.Range("J2").Value = 1741.5555
$ws.Range("K2").Value = 25421.154
$ws.Range("L2").Value = 1741.5555
$ws.Range("M2").Value = -25308.154
$ws.Range("N2").Value = -1967.5555

$ws.Range("H61").Value = 4113.1333
$ws.Range("I61").Value = 1468.5
$ws.Range("K61").Value = 1468.5
$ws.Range("M61").Value = -1256.5

$ws.Range("H74").Value = 2633.75
$ws.Range("I74").Value = 2633.75
$ws.Range("K74").Value = 2633.75
$ws.Range("M74").Value = -1759.75

$ws.Range("H77").Value = 2633.75
$ws.Range("I77").Value = 2633.75
$ws.Range("K77").Value = 13168.75
$ws.Range("M77").Value = -8800.75

$ws.Range("H102").Value = 582.5714
$ws.Range("I102").Value = 582.5714
$ws.Range("K102").Value = 582.5714
$ws.Range("M102").Value = 1039.4286

$ws.Range("H116").Value = 15734.046
$ws.Range("I116").Value = 25421.154
$ws.Range("J116").Value = 1741.5555
$ws.Range("K116").Value = 25421.154
$ws.Range("L116").Value = 1741.5555
$ws.Range("M116").Value = -23127.154
$ws.Range("N116").Value = -6329.5555

$ws.Range("H136").Value = 4113.1333
$ws.Range("I136").Value = 1468.5
$ws.Range("K136").Value = 4405.5
$ws.Range("M136").Value = -1855.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 15734.046
$ws.Range("I3").Value = 25421.154
$ws.Range("J3").Value = 1741.5555
$ws.Range("K3").Value = 25421.154
$ws.Range("L3").Value = 1741.5555
$ws.Range("M3").Value = -25307.154
$ws.Range("N3").Value = -1969.5555

$ws.Range("H99").Value = 19191.652
$ws.Range("I99").Value = 22375.105
$ws.Range("K99").Value = 22375.105
$ws.Range("M99").Value = -20877.105

$ws.Range("H105").Value = 2129.1333
$ws.Range("I105").Value = 2209.7856
$ws.Range("K105").Value = 2209.7856
$ws.Range("M105").Value = -462.7856000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 10027.75
$ws.Range("I94").Value = 17969.166
$ws.Range("J94").Value = 2086.3333
$ws.Range("K94").Value = 17969.166
$ws.Range("L94").Value = 2086.3333
$ws.Range("M94").Value = -17518.166
$ws.Range("N94").Value = -2988.3333

$ws.Range("H105").Value = 692.4091
$ws.Range("I105").Value = 607.4
$ws.Range("J105").Value = 874.5714
$ws.Range("K105").Value = 607.4
$ws.Range("L105").Value = 874.5714
$ws.Range("M105").Value = 1139.6
$ws.Range("N105").Value = -4368.5714

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3717.9375
$ws.Range("J55").Value = 4523.3076
$ws.Range("L55").Value = 13569.9228
$ws.Range("N55").Value = -13923.9228

$ws.Range("H132").Value = 604
$ws.Range("I132").Value = 604
$ws.Range("K132").Value = 5436
$ws.Range("M132").Value = -2906

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 4999.5
$ws.Range("J38").Value = 4999.5
$ws.Range("L38").Value = 4999.5
$ws.Range("N38").Value = -5925.5

$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H63").Value = 44857.145
$ws.Range("J63").Value = 44857.145
$ws.Range("L63").Value = 44857.145
$ws.Range("N63").Value = -46229.145

$ws.Range("H66").Value = 44857.145
$ws.Range("J66").Value = 44857.145
$ws.Range("L66").Value = 134571.435
$ws.Range("N66").Value = -141435.435

$ws.Range("H113").Value = 5402.087
$ws.Range("J113").Value = 10857
$ws.Range("L113").Value = 10857
$ws.Range("N113").Value = -15197

$ws.Range("H122").Value = 3799.64
$ws.Range("I122").Value = 4038.389
$ws.Range("J122").Value = 3185.7144
$ws.Range("K122").Value = 12115.167
$ws.Range("L122").Value = 9557.143199999999
$ws.Range("M122").Value = -9665.167000000001
$ws.Range("N122").Value = -14457.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1862.6428
$ws.Range("I55").Value = 2198.0908
$ws.Range("J55").Value = 632.6667
$ws.Range("K55").Value = 2198.0908
$ws.Range("L55").Value = 632.6667
$ws.Range("M55").Value = -2025.0908
$ws.Range("N55").Value = -978.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 39644.668
$ws.Range("J64").Value = 34467
$ws.Range("L64").Value = 34467
$ws.Range("N64").Value = -34963

$ws.Range("H67").Value = 39644.668
$ws.Range("J67").Value = 34467
$ws.Range("L67").Value = 34467
$ws.Range("N67").Value = -36183
